$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 409
$ws.Range("F5").Value = 8490
$ws.Range("F7").Value = 10574
$ws.Range("F20").Value = 410
$ws.Range("F22").Value = 1803
$ws.Range("F23").Value = 67
$ws.Range("F25").Value = 340
$ws.Range("F26").Value = 283
$ws.Range("F27").Value = 58
$ws.Range("F28").Value = 580
$ws.Range("F30").Value = 1164
$ws.Range("F35").Value = 340
$ws.Range("F38").Value = 128
$ws.Range("F39").Value = 509
$ws.Range("F40").Value = 343
$ws.Range("F42").Value = 280
$ws.Range("F43").Value = 635
$ws.Range("F45").Value = 89
$ws.Range("I42").Value = "//i1.hdslb.com/bfs/openplatform/202409/T5XOZF891727062792168.jpeg"
$ws.Range("I43").Value = "//i1.hdslb.com/bfs/openplatform/202409/rG5Ps2Em1727063078808.jpeg"

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 47
$ws.Range("F17").Value = 380

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 210
$ws.Range("F3").Value = 2796
$ws.Range("F4").Value = 340

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 210
$ws.Range("F6").Value = 340
$ws.Range("F9").Value = 409
$ws.Range("F10").Value = 8490
$ws.Range("F12").Value = 10574
$ws.Range("F19").Value = 1803
$ws.Range("F20").Value = 67
$ws.Range("F22").Value = 283
$ws.Range("F23").Value = 58
$ws.Range("F25").Value = 580
$ws.Range("F28").Value = 1164
$ws.Range("F37").Value = 340
$ws.Range("F38").Value = 128
$ws.Range("F39").Value = 509
$ws.Range("F41").Value = 343
$ws.Range("F43").Value = 280
$ws.Range("F45").Value = 47
$ws.Range("F46").Value = 380
$ws.Range("F47").Value = 635
$ws.Range("F48").Value = 89
$ws.Range("I43").Value = "//i1.hdslb.com/bfs/openplatform/202409/T5XOZF891727062792168.jpeg"
$ws.Range("I47").Value = "//i1.hdslb.com/bfs/openplatform/202409/rG5Ps2Em1727063078808.jpeg"
